$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.827.93'
$ws.Range("E2").Value = '  +1.99%  '
$ws.Range("D3").Value = '3.180.89'
$ws.Range("E3").Value = '  -3.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.78'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '615.38'
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.390'
$ws.Range("E7").Value = '  +1.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.686'
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.998'
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = '3.174.57'
$ws.Range("E10").Value = '  -3.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.572'
$ws.Range("E11").Value = '  -0.99%  '
$ws.Range("E12").Value = '  -6.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  -3.27%  '
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '90.323.29'
$ws.Range("E14").Value = '  +1.33%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.760.21'
$ws.Range("E15").Value = '  -3.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.77'
$ws.Range("E16").Value = '  -4.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.20'
$ws.Range("E17").Value = '  -3.58%  '
$ws.Range("D18").Value = '3.173.88'
$ws.Range("E18").Value = '  -3.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.24'
$ws.Range("E19").Value = '  +3.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.43'
$ws.Range("E20").Value = '  -5.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.25'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000190'
$ws.Range("E22").Value = '  +40.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.53'
$ws.Range("E23").Value = '  -4.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.10'
$ws.Range("E24").Value = '  -5.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.15'
$ws.Range("E25").Value = '  -2.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.83'
$ws.Range("E26").Value = '  -4.42%  '
$ws.Range("D27").Value = '3.440.72'
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '74.69'
$ws.Range("E28").Value = '  -2.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.168'
$ws.Range("E30").Value = '  -9.47%  '
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.16'
$ws.Range("E32").Value = '  +34.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.47'
$ws.Range("E33").Value = '  -4.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '534.23'
$ws.Range("E34").Value = '  -5.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.90'
$ws.Range("E35").Value = '  -4.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.87'
$ws.Range("E36").Value = '  -5.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.25'
$ws.Range("E37").Value = '  -10.77%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '21.96'
$ws.Range("E38").Value = '  -3.30%  '
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.31'
$ws.Range("E39").Value = '  +2.35%  '
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.126'
$ws.Range("E41").Value = '  -9.44%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.376'
$ws.Range("E43").Value = '  -6.41%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.91'
$ws.Range("E44").Value = '  -6.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '144.95'
$ws.Range("E45").Value = '  -5.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '44.69'
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '173.22'
$ws.Range("E47").Value = '  -4.52%  '
$ws.Range("E48").Value = '  -2.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.24'
$ws.Range("E49").Value = '  -5.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.615'
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.08'
$ws.Range("E51").Value = '  -3.91%  '
